$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.943.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3921"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3877"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.360"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08464"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.87"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.200"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.879"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001311"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.654.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06981"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.920"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.008"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.66"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.014.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.039"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.58"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.298"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.829"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.523"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.839.47"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.028"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03020"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08097"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.705"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.87"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2707"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09166"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7554"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.50"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.424"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.31"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6927"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.484"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.091"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08271"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.01"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.402"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.35%  "
